$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.7
$ws.Range("H2").Value = 3.3
$ws.Range("I2").Value = 5.1
$ws.Range("K2").Value = 2.02
$ws.Range("L2").Value = 5.3
$ws.Range("M2").Value = 1.01
$ws.Range("N2").Value = 9.35
$ws.Range("O2").Value = 1.28
$ws.Range("P2").Value = 3.05
$ws.Range("Q2").Value = 1.88
$ws.Range("R2").Value = 1.82
$ws.Range("S2").Value = 1.42
$ws.Range("T2").Value = 2.47
$ws.Range("U2").Value = 1.75
$ws.Range("V2").Value = 1.85
$ws.Range("W2").Value = 6.6
$ws.Range("X2").Value = 7.9
$ws.Range("Y2").Value = 7.8
$ws.Range("Z2").Value = 13.5
$ws.Range("AB2").Value = 25
$ws.Range("AC2").Value = 9.25
$ws.Range("AD2").Value = 6.5
$ws.Range("AE2").Value = 14.5
$ws.Range("AF2").Value = 70
$ws.Range("AH2").Value = 14
$ws.Range("AI2").Value = 32
$ws.Range("AL2").Value = 55
$ws.Range("AM2").Value = 50
$ws.Range("AP2").Value = 17.5
$ws.Range("AR2").Value = 60
$ws.Range("AT2").Value = 2.45
$ws.Range("AV2").Value = 70
$ws.Range("AW2").Value = 6.7
$ws.Range("AX2").Value = 30
$ws.Range("AY2").Value = 35
$ws.Range("AZ2").Value = 200
$ws.Range("BB2").Value = 450
